# [IMP] extras: Improve about page and poetry pdf file.
#
# The "book_name" column (J) contains a small grammar fix: the Spanish
# adjective "Pequeño" (masculine) should agree with the feminine noun
# "colección" and therefore be "Pequeña" (feminine), for both volumes:
#   "Pequeño y humilde colección de poemas Vol. I (Spanish)"
#   "Pequeño y humilde colección de poemas Vol. II (Spanish)"
# become
#   "Pequeña y humilde colección de poemas Vol. I (Spanish)"
#   "Pequeña y humilde colección de poemas Vol. II (Spanish)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poems")

$used = $ws.UsedRange
$used.Replace("Pequeño y humilde colección de poemas", "Pequeña y humilde colección de poemas")
